# Mongolia Summary sheet: the "Source:" block under the sector-distribution
# table is re-laid-out.
#
#   - A blank row is inserted right after the "Source:" caption (row 43),
#     pushing the citation text, the hyperlinked URL, and everything below
#     down by one row.
#   - The hyperlink on the URL cell is removed (it becomes a plain text
#     cell) and the URL text itself is moved one row further down so a
#     blank separator row remains in its old spot.
#   - The long "World Bank (WB), ... Available at ..." citation at the very
#     bottom is shortened down to just "WB".
#
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 43 (inherits the italic "source" style from the
# row it pushes down), shifting rows 43-49 down to 44-50.
$ws.Rows(43).Insert()

# The hyperlinked URL cell is now at A45 - drop the hyperlink itself before
# moving rows around so nothing stale is left referencing the old address.
$ws.Range("A45").Hyperlinks.Delete()

# Remove that row (taking the now plain URL text with it) and immediately
# re-insert a fresh blank row in its place, so row 45 ends up empty again
# and row 46 is free for the URL text (without the hyperlink formatting).
$ws.Rows(45).Delete()
$ws.Rows(45).Insert()
$ws.Range("A46").Value = "http://www.mongolbank.mn/documents/moneypolicy/worldbank/developmentmodule/03.pdf"

# Shorten the final citation line (now at row 50) down to just "WB".
$ws.Range("A50").Value = "WB"
